$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) values are stored as literal text in the workbook even
# though they look numeric. Mark each target cell as Text first so Excel
# does not silently coerce the assigned string into a floating point number.
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
}

Set-TextValue "D2"  "245.75"
Set-TextValue "D3"  "24.19"
Set-TextValue "D4"  "5.329"
Set-TextValue "D5"  "0.05732"
Set-TextValue "D6"  "6.482"
Set-TextValue "D7"  "3.136"
Set-TextValue "D8"  "0.8161"
Set-TextValue "D9"  "0.8697"
Set-TextValue "D10" "0.1378"
Set-TextValue "D11" "0.06994"
Set-TextValue "D12" "0.03178"
Set-TextValue "D13" "0.02913"
Set-TextValue "D14" "0.09400"
Set-TextValue "D15" "3.736"
Set-TextValue "D16" "0.001526"
Set-TextValue "D17" "0.04701"

Set-TextValue "D18" "0.0006009"
$ws.Range("E18").Value = "17OneONE"

Set-TextValue "D19" "0.006185"
Set-TextValue "D20" "0.001239"
Set-TextValue "D21" "0.003869"
Set-TextValue "D22" "0.00008794"
Set-TextValue "D24" "2.150"
Set-TextValue "D25" "0.3175"

Set-TextValue "D40" "0.03718"
Set-TextValue "D41" "0.006436"
Set-TextValue "D42" "0.1057"
Set-TextValue "D43" "0.002274"
Set-TextValue "D44" "0.007847"
Set-TextValue "D45" "0.00005261"

Set-TextValue "D47" "0.3899"

Set-TextValue "D48" "0.004122"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"

$wb.Save()
